$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.725.05"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").Value = "2.289.33"
$ws.Range("E3").Value = "  -1.89%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'96.98"
$ws.Range("E5").Value = "  +1.22%  "

$ws.Range("D6").Value = "'269.71"
$ws.Range("E6").Value = "  -0.45%  "

$ws.Range("E7").Value = "  -0.24%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").Value = "'0.607"
$ws.Range("E9").Value = "  -2.49%  "

$ws.Range("D10").Value = "'45.25"
$ws.Range("E10").Value = "  +0.15%  "

$ws.Range("D11").Value = "'0.0935"
$ws.Range("E11").Value = "  -0.92%  "

$ws.Range("D12").Value = "'7.86"
$ws.Range("E12").Value = "  -2.67%  "

$ws.Range("D13").Value = "'0.106"
$ws.Range("E13").Value = "  +1.23%  "

$ws.Range("D14").Value = "'15.91"
$ws.Range("E14").Value = "  +0.82%  "

$ws.Range("D15").Value = "2.632.34"
$ws.Range("E15").Value = "  -1.72%  "

$ws.Range("D16").Value = "'0.858"
$ws.Range("E16").Value = "  -0.15%  "

$ws.Range("D17").Value = "2.299.13"
$ws.Range("E17").Value = "  -1.27%  "

$ws.Range("D18").Value = "43.728.14"
$ws.Range("E18").Value = "  +0.11%  "

$ws.Range("E19").Value = "  +1.69%  "

$ws.Range("E20").Value = "  -2.95%  "

$ws.Range("D21").Value = "'72.12"
$ws.Range("E21").Value = "  +0.28%  "

$ws.Range("D22").Value = "'2.46"
$ws.Range("E22").Value = "  +8.92%  "

$ws.Range("D23").Value = "'232.79"
$ws.Range("E23").Value = "  -4.95%  "

$ws.Range("D24").Value = "'9.04"
$ws.Range("E24").Value = "  -3.92%  "

$ws.Range("D25").Value = "'2.71"
$ws.Range("E25").Value = "  +6.45%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").Value = "'11.23"
$ws.Range("E27").Value = "  -1.53%  "

$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("D29").Value = "'2.29"
$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("D30").Value = "'38.55"
$ws.Range("E30").Value = "  +0.32%  "

$ws.Range("D31").Value = "'176.30"
$ws.Range("E31").Value = "  +1.98%  "

$ws.Range("D32").Value = "'21.80"
$ws.Range("E32").Value = "  -3.70%  "

$ws.Range("D33").Value = "'0.0893"
$ws.Range("E33").Value = "  -0.43%  "

$ws.Range("E34").Value = "  -1.69%  "

$ws.Range("E35").Value = "  +0.82%  "

$ws.Range("D36").Value = "'4.67"
$ws.Range("E36").Value = "  +7.03%  "

$ws.Range("D37").Value = "'0.108"
$ws.Range("E37").Value = "  +0.53%  "

$ws.Range("E38").Value = "  -2.33%  "

$ws.Range("D39").Value = "'3.50"
$ws.Range("E39").Value = "  +4.00%  "

$ws.Range("D40").Value = "'0.237"
$ws.Range("E40").Value = "  +1.42%  "

$ws.Range("E41").Value = "  -2.55%  "

$ws.Range("D42").Value = "'12.20"
$ws.Range("E42").Value = "  +0.63%  "

$ws.Range("E43").Value = "  -2.27%  "

$ws.Range("D44").Value = "'64.32"
$ws.Range("E44").Value = "  +3.90%  "

$ws.Range("D45").Value = "'5.21"
$ws.Range("E45").Value = "  -3.13%  "

$ws.Range("E46").Value = "  -4.49%  "

$ws.Range("E47").Value = "  -0.74%  "

$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").Value = "'1.23"
$ws.Range("E48").Value = "  +1.15%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'98.93"
$ws.Range("E49").Value = "  -1.75%  "

$ws.Range("D50").Value = "'0.444"
$ws.Range("E50").Value = "  +7.27%  "

$ws.Range("D51").Value = "'1.51"
$ws.Range("E51").Value = "  +11.06%  "
